# Insert a new data row at row 124, shifting the existing rows 124-182
# down to 125-183 (the last row's data moves to the new row 183), then
# populate the freshly inserted row 124 with its new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value = 9
$ws.Cells.Item(124, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(124, 3).Value = "Metropolitana"
$ws.Cells.Item(124, 4).Value = 44452
$ws.Cells.Item(124, 5).Value = 13
$ws.Cells.Item(124, 6).Value = 100112052
$ws.Cells.Item(124, 7).Value = "Albahaca"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 196
$ws.Cells.Item(124, 11).Value = 4500
$ws.Cells.Item(124, 12).Value = 5000
$ws.Cells.Item(124, 13).Value = 4750
$ws.Cells.Item(124, 14).Value = "$/paquete"
$ws.Cells.Item(124, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(124, 16).Value = 4750
$ws.Cells.Item(124, 17).Value = 1
$ws.Cells.Item(124, 18).Value = "Hortaliza"
